$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1551.3903
$ws.Range("J17").Value = 1552.925
$ws.Range("L17").Value = 4658.775
$ws.Range("N17").Value = -4994.775
$ws.Range("H33").Value = 196
$ws.Range("I33").Value = 196
$ws.Range("K33").Value = 196
$ws.Range("M33").Value = 33
$ws.Range("H42").Value = 350
$ws.Range("I42").Value = 200
$ws.Range("K42").Value = 600
$ws.Range("M42").Value = -370
$ws.Range("H74").Value = 12830.286
$ws.Range("I74").Value = 13128.263
$ws.Range("K74").Value = 13128.263
$ws.Range("M74").Value = -12192.263
$ws.Range("H77").Value = 12830.286
$ws.Range("I77").Value = 13128.263
$ws.Range("K77").Value = 65641.315
$ws.Range("M77").Value = -60961.315
$ws.Range("H80").Value = 1632.125
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = $null
$ws.Range("N81").Value = 0
$ws.Range("H83").Value = 1632.125
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = $null
$ws.Range("N84").Value = 0
$ws.Range("H98").Value = 3636.8096
$ws.Range("J98").Value = 7800
$ws.Range("L98").Value = 7800
$ws.Range("N98").Value = -10796
$ws.Range("H113").Value = 3965.3
$ws.Range("I113").Value = 3428.6
$ws.Range("J113").Value = 4502
$ws.Range("K113").Value = 3428.6
$ws.Range("L113").Value = 4502
$ws.Range("M113").Value = -174.5999999999999
$ws.Range("N113").Value = -11010
$ws.Range("H122").Value = 3636.8096
$ws.Range("J122").Value = 7800
$ws.Range("L122").Value = 23400
$ws.Range("N122").Value = -28300
$ws.Range("H132").Value = 5735
$ws.Range("I132").Value = 6085.8
$ws.Range("K132").Value = 18257.4
$ws.Range("M132").Value = -15727.4
$ws.Range("H137").Value = 1788342.4
$ws.Range("I137").Value = 2633001.2
$ws.Range("J137").Value = 5173.4443
$ws.Range("K137").Value = 7899003.600000001
$ws.Range("L137").Value = 15520.3329
$ws.Range("M137").Value = -7896453.600000001
$ws.Range("N137").Value = -20620.3329
$ws.Range("H138").Value = 4326.517
$ws.Range("I138").Value = 18746
$ws.Range("J138").Value = 3258.4075
$ws.Range("K138").Value = 56238
$ws.Range("L138").Value = 9775.2225
$ws.Range("M138").Value = -51098
$ws.Range("N138").Value = -20055.2225

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 11000
$ws.Range("I19").Value = 17000
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 17000
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = -16771
$ws.Range("N19").Value = -5458
$ws.Range("H40").Value = 23014
$ws.Range("I40").Value = 16028
$ws.Range("K40").Value = 16028
$ws.Range("M40").Value = -15852
$ws.Range("H50").Value = 3381.6667
$ws.Range("J50").Value = 3381.6667
$ws.Range("L50").Value = 3381.6667
$ws.Range("N50").Value = -4809.6667
$ws.Range("H61").Value = 2806.2
$ws.Range("I61").Value = 1832.7142
$ws.Range("K61").Value = 1832.7142
$ws.Range("M61").Value = -1620.7142
$ws.Range("H74").Value = 182280.3
$ws.Range("I74").Value = 294881.75
$ws.Range("K74").Value = 294881.75
$ws.Range("M74").Value = -294007.75
$ws.Range("H77").Value = 182280.3
$ws.Range("I77").Value = 294881.75
$ws.Range("K77").Value = 1474408.75
$ws.Range("M77").Value = -1470040.75
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = $null
$ws.Range("N92").Value = 0
$ws.Range("H122").Value = 2394.3235
$ws.Range("I122").Value = 2387
$ws.Range("K122").Value = 7161
$ws.Range("M122").Value = -4711
$ws.Range("H132").Value = 2113.7908
$ws.Range("I132").Value = 1961.8214
$ws.Range("J132").Value = 2397.4666
$ws.Range("K132").Value = 5885.4642
$ws.Range("L132").Value = 7192.399800000001
$ws.Range("M132").Value = -3355.4642
$ws.Range("N132").Value = -12252.3998
$ws.Range("H136").Value = 2806.2
$ws.Range("I136").Value = 1832.7142
$ws.Range("K136").Value = 5498.142599999999
$ws.Range("M136").Value = -2948.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 16251803
$ws.Range("I105").Value = 1112712
$ws.Range("K105").Value = 1112712
$ws.Range("M105").Value = -1110965
$ws.Range("H134").Value = 1981.091
$ws.Range("I134").Value = 1524.4857
$ws.Range("K134").Value = 4573.4571
$ws.Range("M134").Value = -2038.4571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 56641.125
$ws.Range("J9").Value = 56641.125
$ws.Range("L9").Value = 56641.125
$ws.Range("N9").Value = -56977.125
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10574
$ws.Range("H31").Value = 6043.0293
$ws.Range("I31").Value = 4630.364
$ws.Range("J31").Value = 6718.6523
$ws.Range("K31").Value = 4630.364
$ws.Range("L31").Value = 6718.6523
$ws.Range("M31").Value = -4335.364
$ws.Range("N31").Value = -7308.6523
$ws.Range("H34").Value = 6043.0293
$ws.Range("I34").Value = 4630.364
$ws.Range("J34").Value = 6718.6523
$ws.Range("K34").Value = 4630.364
$ws.Range("L34").Value = 6718.6523
$ws.Range("M34").Value = -4428.364
$ws.Range("N34").Value = -7122.6523
$ws.Range("H58").Value = 3879.8
$ws.Range("J58").Value = 4110.8887
$ws.Range("L58").Value = 4110.8887
$ws.Range("N58").Value = -4516.8887
$ws.Range("H132").Value = 15878474
$ws.Range("I132").Value = 3661.923
$ws.Range("K132").Value = 10985.769
$ws.Range("M132").Value = -8455.769
$ws.Range("H136").Value = 3879.8
$ws.Range("J136").Value = 4110.8887
$ws.Range("L136").Value = 12332.6661
$ws.Range("N136").Value = -17432.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 430
$ws.Range("I14").Value = 430
$ws.Range("K14").Value = 1290
$ws.Range("M14").Value = -1117
$ws.Range("H17").Value = 1050
$ws.Range("I17").Value = 1050
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 3150
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = -2981
$ws.Range("H82").Value = 10749.5
$ws.Range("H85").Value = 10749.5
$ws.Range("H109").Value = 3300.3333
$ws.Range("I109").Value = 709.5
$ws.Range("K109").Value = 2128.5
$ws.Range("M109").Value = -1088.5
$ws.Range("H124").Value = 5500
$ws.Range("I124").Value = 5500
$ws.Range("K124").Value = 16500
$ws.Range("M124").Value = -11590

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = $null
$ws.Range("N52").Value = 0
$ws.Range("H92").Value = 9999
$ws.Range("J92").Value = 9999
$ws.Range("L92").Value = 9999
$ws.Range("N92").Value = -13743
$ws.Range("H102").Value = 4487.9062
$ws.Range("I102").Value = 1204.3478
$ws.Range("J102").Value = 12879.223
$ws.Range("K102").Value = 1204.3478
$ws.Range("L102").Value = 12879.223
$ws.Range("M102").Value = 417.6522
$ws.Range("N102").Value = -16123.223
$ws.Range("H126").Value = 7288
$ws.Range("I126").Value = 2181
$ws.Range("K126").Value = 6543
$ws.Range("M126").Value = -4073

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3366.5417
$ws.Range("I7").Value = 2540.35
$ws.Range("K7").Value = 2540.35
$ws.Range("M7").Value = -2428.35
$ws.Range("H16").Value = 1998.125
$ws.Range("I16").Value = 1140.7142
$ws.Range("K16").Value = 1140.7142
$ws.Range("M16").Value = -970.7141999999999
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = $null
$ws.Range("M25").Value = 10000
$ws.Range("N25").Value = -10460
$ws.Range("H43").Value = 10430
$ws.Range("I43").Value = 6670.6665
$ws.Range("K43").Value = 6670.6665
$ws.Range("M43").Value = -6477.6665
$ws.Range("H46").Value = 3641.625
$ws.Range("I46").Value = 2751.4546
$ws.Range("K46").Value = 2751.4546
$ws.Range("M46").Value = -2563.4546
$ws.Range("H126").Value = 3366.5417
$ws.Range("I126").Value = 2540.35
$ws.Range("K126").Value = 7621.049999999999
$ws.Range("M126").Value = -5151.049999999999
$ws.Range("H132").Value = 4624.1284
$ws.Range("I132").Value = 3110.077
$ws.Range("J132").Value = 7652.231
$ws.Range("K132").Value = 9330.231
$ws.Range("L132").Value = 22956.693
$ws.Range("M132").Value = -6800.231
$ws.Range("N132").Value = -28016.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("H113").Value = 928.8889
$ws.Range("J113").Value = 1197
$ws.Range("L113").Value = 3591
$ws.Range("N113").Value = -7931
